$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Tick the four checkboxes that were unchecked before and are checked now.
#    Each Word content-control checkbox stores its state twice: the SDT's
#    w14:checked flag, and the glyph in the run inside sdtContent (BALLOT BOX
#    U+2610 "☐" <-> BALLOT BOX WITH X U+2612 "☒"). Both need updating.
# ---------------------------------------------------------------------------
$checkboxIds = @(-2098477851, 98995749, 993686833, -1437672738, -495417201)

$ccs = $d.ContentControls
for ($i = 1; $i -le $ccs.Count; $i++) {
    $cc = $ccs.Item($i)
    foreach ($id in $checkboxIds) {
        if ($cc.ID -eq $id) {
            $cc.Checked = $true
            $cc.Range.Text = [char]0x2612
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Move the stray "_GoBack" bookmark: delete it from after "Betalingsstatus"
#    and re-create it (collapsed) right after "Brainstorm: problemformulering".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$findRng = $d.Content
$found = $findRng.Find.Execute("Brainstorm: problemformulering", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $findRng.End

    # Placing a bookmark exactly at a paragraph-end offset directly can misfire,
    # so nudge the boundary out of the way first: append a throwaway run after
    # the target text, drop the bookmark in the now-safe interior position
    # between the two runs, then remove the throwaway run again.
    $findRng.InsertAfter("x")

    $bmRng = $d.Range($insertPoint, $insertPoint)
    $d.Bookmarks.Add("_GoBack", $bmRng)

    $cleanupRng = $d.Range($insertPoint, $insertPoint + 1)
    $cleanupRng.Text = ""
}
